$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for Chirimoya (Macroferia Regional de Talca) needs
# to be inserted as the new row 6; the previous row 6 (and everything below
# it) shifts down by one row. Insert a blank row at position 7 so the
# existing row 6 stays put and rows 7-16 move to 8-17.
$ws.Rows.Item(7).Insert()

# The row that used to be row 6 is now duplicated down into the freshly
# inserted row 7 (all columns A:T), preserving that record unchanged.
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(7, $col).Value = $ws.Cells.Item(6, $col).Value2
}

# Row 6 becomes this week's new entry: same market/product/quality, but a
# newer date and an updated Volumen figure.
$ws.Cells.Item(6, 4).Value = 44459   # D6 - Fecha
$ws.Cells.Item(6, 13).Value = 100    # M6 - Volumen
